$d = $word.ActiveDocument

# ------------------------------------------------------------------
# REPORTGEN-397 : update extension name for quality standards
#
# The footnote sentence that starts with "**" currently reads:
#   ...extension "Standard Quality Rules" is installed...
# and must become:
#   ...extension "Quality Standards Support" is installed...
#
# In addition, the stray "_GoBack" bookmark that currently sits
# (collapsed) in front of the second, non-"**" copy of this same
# sentence further down needs to move so that it wraps the newly
# inserted "Quality Standards Support" text instead.
# ------------------------------------------------------------------

# 1) Locate the "**" sentence (only that copy of the sentence is edited).
$sentence = $d.Content
$sentence.Find.ClearFormatting()
$sentenceFound = $sentence.Find.Execute( `
    "** The selection of metrics by standard quality tag name", `
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $sentenceFound) {
    throw "Could not find the '**' sentence to edit"
}
$sentenceStart = $sentence.Start

# 2) Within (and after) that sentence, locate the quoted extension name.
$scope = $d.Range($sentenceStart, $d.Content.End)
$scope.Find.ClearFormatting()
$quotedFound = $scope.Find.Execute("Standard Quality Rules", $true, $false, `
    $false, $false, $false, $true, 1, $false, "", 0)
if (-not $quotedFound) {
    throw "Could not find the quoted extension name"
}
$oldTextStart = $scope.Start
$oldTextEnd = $scope.End

$newName = "Quality Standards Support"

# 3) Replace the quoted extension name with the new one. At this point
#    the edit still lives inside the single original run, so Word keeps
#    (and correctly re-evaluates) the run's own xml:space handling.
$editRange = $d.Range($oldTextStart, $oldTextEnd)
$editRange.Text = $newName
$newTextEnd = $oldTextStart + $newName.Length

# 4) Force a clean run split exactly around the new text by toggling a
#    character property (and immediately reverting it) on that precise
#    sub-range. This carves the paragraph into three runs - before /
#    new-name / after - without touching their text content again.
$newNameRange = $d.Range($oldTextStart, $newTextEnd)
$newNameRange.Bold = 1
$newNameRange.Bold = 0

# 5) Wrap the (now isolated) new-name run with the "_GoBack" bookmark.
#    Re-using the existing bookmark name automatically relocates it,
#    removing the stale bookmarkStart/bookmarkEnd pair that used to sit
#    in front of the duplicate sentence further down the document.
$bookmarkRange = $d.Range($oldTextStart, $newTextEnd)
$d.Bookmarks.Add("_GoBack", $bookmarkRange)
